# Fruta / hortaliza, semanal
# Insert a new weekly record for "Santina" / "Primera" (Provincia de Curicó)
# as row 40, pushing the existing Cereza records (old rows 40-53) down by
# one row to rows 41-54.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(40).Insert()

$ws.Cells.Item(40, 1).Value  = 11
$ws.Cells.Item(40, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(40, 3).Value  = "Bíobío"
$ws.Cells.Item(40, 4).Value  = 44523
$ws.Cells.Item(40, 5).Value  = 8
$ws.Cells.Item(40, 6).Value  = "Fruta"
$ws.Cells.Item(40, 7).Value  = 100103
$ws.Cells.Item(40, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(40, 9).Value  = 100103001
$ws.Cells.Item(40, 10).Value = "Cereza"
$ws.Cells.Item(40, 11).Value = "Santina"
$ws.Cells.Item(40, 12).Value = "Primera"
$ws.Cells.Item(40, 13).Value = 70
$ws.Cells.Item(40, 14).Value = 20000
$ws.Cells.Item(40, 15).Value = 22000
$ws.Cells.Item(40, 16).Value = 20571
$ws.Cells.Item(40, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(40, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(40, 19).Value = 2057
$ws.Cells.Item(40, 20).Value = 10
